# Update the handback/handoff timestamps in the zh-cn and de-de sheets
# to reflect the newly generated report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 5 corresponds to the 7c3ff293-... handback entry
$wsZhCn.Range("D5").Value = "2016-02-18 10:24:07"
$wsZhCn.Range("G5").Value = "2016-02-18 10:24:57"

# de-de sheet: row 5 corresponds to the 7c3ff293-... handback entry
$wsDeDe.Range("D5").Value = "2016-02-18 10:24:18"
$wsDeDe.Range("G5").Value = "2016-02-18 10:25:20"
